$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.43424333333333
$ws.Range("H2").Value = 31.30273
$ws.Range("I2").Value = 0.9711091978791583
$ws.Range("J2").Value = 0.9711091978791584
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6559766666666667
$ws.Range("N2").Value = 1.96793
$ws.Range("O2").Value = 0.030799191223283
$ws.Range("P2").Value = 0.030799191223283
$ws.Range("Q2").Value = 6.844620160988889
$ws.Range("R2").Value = 61.60158144889999
$ws.Range("S2").Value = 0.02990937788416917
$ws.Range("T2").Value = 0.02990937788416917

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.43424333333333
$ws.Range("H3").Value = 31.30273
$ws.Range("I3").Value = 0.9711091978791583
$ws.Range("J3").Value = 0.9711091978791584
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 16.81477433333333
$ws.Range("N3").Value = 50.444323
$ws.Range("O3").Value = 0.7894815111340611
$ws.Range("P3").Value = 0.789481511134061
$ws.Range("Q3").Value = 175.4494469890878
$ws.Range("R3").Value = 1579.04502290179
$ws.Range("S3").Value = 0.7666727570178238
$ws.Range("T3").Value = 0.7666727570178238

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.43424333333333
$ws.Range("H4").Value = 31.30273
$ws.Range("I4").Value = 0.9711091978791583
$ws.Range("J4").Value = 0.9711091978791584
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5954103333333334
$ws.Range("N4").Value = 1.786231
$ws.Range("O4").Value = 0.02795550153610953
$ws.Range("P4").Value = 0.02795550153610953
$ws.Range("Q4").Value = 6.212656301181111
$ws.Range("R4").Value = 55.91390671062999
$ws.Range("S4").Value = 0.0271478446730409
$ws.Range("T4").Value = 0.0271478446730409

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.43424333333333
$ws.Range("H5").Value = 31.30273
$ws.Range("I5").Value = 0.9711091978791583
$ws.Range("J5").Value = 0.9711091978791584
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.232341666666667
$ws.Range("N5").Value = 9.697025
$ws.Range("O5").Value = 0.1517637961065464
$ws.Range("P5").Value = 0.1517637961065464
$ws.Range("Q5").Value = 33.72703948647222
$ws.Range("R5").Value = 303.5433553782499
$ws.Range("S5").Value = 0.1473792183041244
$ws.Range("T5").Value = 0.1473792183041244

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.310422
$ws.Range("H6").Value = 0.9312659999999999
$ws.Range("I6").Value = 0.02889080212084161
$ws.Range("J6").Value = 0.02889080212084161
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6559766666666667
$ws.Range("N6").Value = 1.96793
$ws.Range("O6").Value = 0.030799191223283
$ws.Range("P6").Value = 0.030799191223283
$ws.Range("Q6").Value = 0.20362958882
$ws.Range("R6").Value = 1.83266629938
$ws.Range("S6").Value = 0.0008898133391138307
$ws.Range("T6").Value = 0.0008898133391138308

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.310422
$ws.Range("H7").Value = 0.9312659999999999
$ws.Range("I7").Value = 0.02889080212084161
$ws.Range("J7").Value = 0.02889080212084161
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 16.81477433333333
$ws.Range("N7").Value = 50.444323
$ws.Range("O7").Value = 0.7894815111340611
$ws.Range("P7").Value = 0.789481511134061
$ws.Range("Q7").Value = 5.219675878102
$ws.Range("R7").Value = 46.97708290291799
$ws.Range("S7").Value = 0.02280875411623717
$ws.Range("T7").Value = 0.02280875411623717

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.310422
$ws.Range("H8").Value = 0.9312659999999999
$ws.Range("I8").Value = 0.02889080212084161
$ws.Range("J8").Value = 0.02889080212084161
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5954103333333334
$ws.Range("N8").Value = 1.786231
$ws.Range("O8").Value = 0.02795550153610953
$ws.Range("P8").Value = 0.02795550153610953
$ws.Range("Q8").Value = 0.184828466494
$ws.Range("R8").Value = 1.663456198446
$ws.Range("S8").Value = 0.0008076568630686239
$ws.Range("T8").Value = 0.000807656863068624

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.310422
$ws.Range("H9").Value = 0.9312659999999999
$ws.Range("I9").Value = 0.02889080212084161
$ws.Range("J9").Value = 0.02889080212084161
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.232341666666667
$ws.Range("N9").Value = 9.697025
$ws.Range("O9").Value = 0.1517637961065464
$ws.Range("P9").Value = 0.1517637961065464
$ws.Range("Q9").Value = 1.00338996485
$ws.Range("R9").Value = 9.03050968365
$ws.Range("S9").Value = 0.004384577802421984
$ws.Range("T9").Value = 0.004384577802421985
